# Finished Week 13 logging
# Update Target Depth Data (Row 2 == "H") for both OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 297
$wsOff.Range("C2").Value = 198
$wsOff.Range("D2").Value = 139
$wsOff.Range("E2").Value = 62
$wsOff.Range("G2").Value = 11

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 434
$wsDef.Range("C2").Value = 328
$wsDef.Range("D2").Value = 85
$wsDef.Range("E2").Value = 52
